$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaIndexByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# --- Change 1: remove the "Idle" process_list entry paragraph entirely ---
$idleIdx = Find-ParaIndexByPrefix "`tcomment: Test"
if ($idleIdx -gt 0) {
    $d.Paragraphs($idleIdx).Range.Delete() | Out-Null
}

# --- Change 2: insert a new kernel_list entry (lltdio.sys) right after the
#     TSDDD.dll entry paragraph, before the "network_list" heading ---
$tsdIdx = Find-ParaIndexByPrefix "`tcomment: Good"
$tsdRange = $d.Paragraphs($tsdIdx).Range
$tsdRange.InsertParagraphAfter() | Out-Null

$newIdx = $tsdIdx + 1
$newRange = $d.Paragraphs($newIdx).Range
$lltdioXml = "<w:p $wns><w:r><w:tab/><w:t>comment: dea</w:t><w:br/><w:tab/><w:t>name: lltdio.sys</w:t><w:br/><w:tab/><w:t>physical_offset: 2107855200</w:t><w:br/><w:tab/><w:t>marked: disabled</w:t><w:br/><w:tab/><w:t>object_id: 128</w:t><w:br/><w:tab/><w:t>file_path: /SystemRoot/system32/DRIVERS/lltdio.sys</w:t><w:br/></w:r></w:p>"
$newRange.InsertXML($lltdioXml) | Out-Null

# --- Change 3: reorder lines inside the network_list (UDPv6) entry paragraph ---
$netIdx = Find-ParaIndexByPrefix "`tcomment: udl"
$netRange = $d.Paragraphs($netIdx).Range
$netXml = "<w:p $wns><w:r><w:tab/><w:t>comment: udl</w:t><w:br/><w:tab/><w:t>pid: 724</w:t><w:br/><w:tab/><w:t>physical_offset: 2107714752</w:t><w:br/><w:tab/><w:t>object_id: 607</w:t><w:br/><w:tab/><w:t>marked: disabled</w:t><w:br/><w:tab/><w:t>local_address: 65152:0:55640:35652:22209:54695</w:t><w:br/><w:tab/><w:t>protocol_version: UDPv6</w:t><w:br/><w:tab/><w:t>owner_process: svchost.exe</w:t><w:br/><w:tab/><w:t>port: 1900</w:t><w:br/></w:r></w:p>"
$netRange.InsertXML($netXml) | Out-Null
